# Swap the presentation's theme palette: the deck currently uses the
# "Integral" theme (ppt/theme/theme2.xml, bound to the slide master) and
# carries an unused "Office Theme" theme part (ppt/theme/theme1.xml,
# bound only to the notes master). The commit swaps the two themes'
# content so the slide master ends up on "Office Theme" colors.
#
# The PowerPoint colour-scheme object model only exposes the *active*
# (slide-master) theme for editing, via Slide.ThemeColorScheme, so we
# rewrite every slot of that live theme to the "Office Theme" palette
# (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink). dk1 and lt1 are identical
# between the two themes (black/white); the remaining ten swap over.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette = the "Office Theme" colours (RGB ints, 0x00BBGGRR order)
$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
